$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: relocate "months_as_a_registered" so it sits right after
# "bd" (it currently lives as the very last row of the sheet).
# ------------------------------------------------------------------
$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$monthsRow = -1
$bdRow = -1
for ($i = 1; $i -le $rowCount; $i++) {
    $lbl = $ws.Cells.Item($i, 1).Value()
    if ($lbl -eq "months_as_a_registered") { $monthsRow = $i }
    if ($lbl -eq "bd") { $bdRow = $i }
}

if ($monthsRow -ne -1 -and $bdRow -ne -1) {
    $labelVal = $ws.Cells.Item($monthsRow, 1).Value()
    $numVal = $ws.Cells.Item($monthsRow, 2).Value()

    # remove it from its old position
    $ws.Rows.Item($monthsRow).Delete()

    # insert it just below "bd"
    $insertAt = $bdRow + 1
    $ws.Rows.Item($insertAt).Insert()
    $ws.Cells.Item($insertAt, 1).Value = $labelVal
    $ws.Cells.Item($insertAt, 2).Value = $numVal
}

# ------------------------------------------------------------------
# Step 2: drop the "*_mov_..._m3" (3-month moving window) metrics and
# the "total_secs" row entirely, keeping the "_m6" variants.
# ------------------------------------------------------------------
$used2 = $ws.UsedRange
$rowCount2 = $used2.Rows.Count
for ($i = $rowCount2; $i -ge 1; $i--) {
    $lbl2 = $ws.Cells.Item($i, 1).Value()
    if ($lbl2 -match "_mov_.*_m3$" -or $lbl2 -eq "total_secs") {
        $ws.Rows.Item($i).Delete()
    }
}
